$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Estimated Hours (column B)
$ws.Range("B4").Value = 2
$ws.Range("B6").Value = 3
$ws.Range("B7").Value = 3
$ws.Range("B8").Value = 1
$ws.Range("B9").Value = 2

# Update Actual Hours (column C)
$ws.Range("C3").Value = 0
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 2
$ws.Range("C6").Value = 2
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 7.5
$ws.Range("C10").Value = 8

# Add totals row
$ws.Range("B11").Formula = "=SUM(B3:B10)"
$ws.Range("C11").Formula = "=SUM(C3:C10)"

# Update selection to reflect last-used cell
$ws.Range("C10").Select()
